# Generate Report for Handback
#
# Refreshes the "Latest Handback DateTime" (column L) for the first data
# row (row 2, file 28c8dcfe-f72d-4a1e-8572-60038800e9d0...) on both the
# "zh-cn" and "de-de" report sheets, reflecting a newly generated handback
# report timestamp.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("L2").Value = "2016-12-16 09:39:12"
$wsDe.Range("L2").Value = "2016-12-16 09:39:30"
